$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1074.5
$ws.Cells.Item(18, 9).Value = 1074.5
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 1074.5
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = -790.5
$ws.Cells.Item(18, 14).ClearContents()
$ws.Cells.Item(98, 8).Value = 2160.2666
$ws.Cells.Item(98, 9).Value = 634.3333
$ws.Cells.Item(98, 10).Value = 3177.5557
$ws.Cells.Item(98, 11).Value = 634.3333
$ws.Cells.Item(98, 12).Value = 3177.5557
$ws.Cells.Item(98, 13).Value = 863.6667
$ws.Cells.Item(122, 8).Value = 2160.2666
$ws.Cells.Item(122, 9).Value = 634.3333
$ws.Cells.Item(122, 10).Value = 3177.5557
$ws.Cells.Item(122, 11).Value = 1902.9999
$ws.Cells.Item(122, 12).Value = 9532.667099999999
$ws.Cells.Item(122, 13).Value = 547.0001
$ws.Cells.Item(131, 8).Value = 6745
$ws.Cells.Item(131, 9).Value = 1353.8889
$ws.Cells.Item(131, 10).Value = 18875
$ws.Cells.Item(131, 11).Value = 4061.6667
$ws.Cells.Item(131, 12).Value = 56625
$ws.Cells.Item(131, 13).Value = 978.3333000000002
$ws.Cells.Item(132, 8).Value = 435.4
$ws.Cells.Item(132, 9).Value = 455.7857
$ws.Cells.Item(132, 10).Value = 150
$ws.Cells.Item(132, 11).Value = 1367.3571
$ws.Cells.Item(132, 12).Value = 450
$ws.Cells.Item(132, 13).Value = 1162.6429

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 896
$ws.Cells.Item(2, 9).Value = 844
$ws.Cells.Item(2, 10).Value = 1000
$ws.Cells.Item(2, 11).Value = 844
$ws.Cells.Item(2, 12).Value = 1000
$ws.Cells.Item(2, 13).Value = -731
$ws.Cells.Item(2, 14).Value = -1226
$ws.Cells.Item(61, 8).Value = 2366.9473
$ws.Cells.Item(61, 9).Value = 2366.9473
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 2366.9473
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -2154.9473
$ws.Cells.Item(110, 8).Value = 12333800
$ws.Cells.Item(110, 9).Value = 12333800
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 11).Value = 12333800
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 13).Value = -12331755
$ws.Cells.Item(116, 8).Value = 896
$ws.Cells.Item(116, 9).Value = 844
$ws.Cells.Item(116, 10).Value = 1000
$ws.Cells.Item(116, 11).Value = 844
$ws.Cells.Item(116, 12).Value = 1000
$ws.Cells.Item(116, 13).Value = 1450
$ws.Cells.Item(116, 14).Value = -5588
$ws.Cells.Item(136, 8).Value = 2366.9473
$ws.Cells.Item(136, 9).Value = 2366.9473
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 7100.841899999999
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -4550.841899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 896
$ws.Cells.Item(3, 9).Value = 844
$ws.Cells.Item(3, 10).Value = 1000
$ws.Cells.Item(3, 11).Value = 844
$ws.Cells.Item(3, 12).Value = 1000
$ws.Cells.Item(3, 13).Value = -730
$ws.Cells.Item(3, 14).Value = -1228
$ws.Cells.Item(42, 8).Value = 495000
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 495000
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 495000
$ws.Cells.Item(42, 14).Value = -495656
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).ClearContents()
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 13).ClearContents()
$ws.Cells.Item(80, 8).Value = 602.8823
$ws.Cells.Item(80, 9).Value = 716.75
$ws.Cells.Item(80, 10).Value = 329.6
$ws.Cells.Item(80, 11).Value = 716.75
$ws.Cells.Item(80, 12).Value = 329.6
$ws.Cells.Item(80, 13).Value = 281.25
$ws.Cells.Item(83, 8).Value = 602.8823
$ws.Cells.Item(83, 9).Value = 716.75
$ws.Cells.Item(83, 10).Value = 329.6
$ws.Cells.Item(83, 11).Value = 3583.75
$ws.Cells.Item(83, 12).Value = 1648
$ws.Cells.Item(83, 13).Value = 1408.25
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents()
$ws.Cells.Item(94, 8).Value = 963.1429000000001
$ws.Cells.Item(94, 9).Value = 963.1429000000001
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 963.1429000000001
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = -512.1429000000001
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 11).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 14).ClearContents()
$ws.Cells.Item(99, 8).Value = 3111
$ws.Cells.Item(99, 9).Value = 2222
$ws.Cells.Item(99, 10).Value = 4000
$ws.Cells.Item(99, 11).Value = 2222
$ws.Cells.Item(99, 12).Value = 4000
$ws.Cells.Item(99, 13).Value = -724
$ws.Cells.Item(99, 14).Value = -6996
$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 3781.3572
$ws.Cells.Item(107, 9).Value = 3403.5454
$ws.Cells.Item(107, 10).Value = 5166.6665
$ws.Cells.Item(107, 11).Value = 3403.5454
$ws.Cells.Item(107, 12).Value = 5166.6665
$ws.Cells.Item(107, 13).Value = -1483.5454
$ws.Cells.Item(107, 14).Value = -9006.666499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2309.5386
$ws.Cells.Item(31, 9).Value = 1930.125
$ws.Cells.Item(31, 10).Value = 2916.6
$ws.Cells.Item(31, 11).Value = 1930.125
$ws.Cells.Item(31, 12).Value = 2916.6
$ws.Cells.Item(31, 13).Value = -1635.125
$ws.Cells.Item(34, 8).Value = 2309.5386
$ws.Cells.Item(34, 9).Value = 1930.125
$ws.Cells.Item(34, 10).Value = 2916.6
$ws.Cells.Item(34, 11).Value = 1930.125
$ws.Cells.Item(34, 12).Value = 2916.6
$ws.Cells.Item(34, 13).Value = -1728.125
$ws.Cells.Item(99, 8).Value = 3713.8572
$ws.Cells.Item(99, 9).Value = 1999
$ws.Cells.Item(99, 10).Value = 3999.6667
$ws.Cells.Item(99, 11).Value = 1999
$ws.Cells.Item(99, 12).Value = 3999.6667
$ws.Cells.Item(99, 13).Value = -501
$ws.Cells.Item(99, 14).Value = -6995.6667
$ws.Cells.Item(108, 8).Value = 24999.5
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 24999.5
$ws.Cells.Item(108, 11).Value = 0
$ws.Cells.Item(108, 12).Value = 24999.5
$ws.Cells.Item(108, 14).Value = -32679.5
$ws.Cells.Item(122, 8).Value = 1374.8125
$ws.Cells.Item(122, 9).Value = 1398.8667
$ws.Cells.Item(122, 10).Value = 1014
$ws.Cells.Item(122, 11).Value = 4196.6001
$ws.Cells.Item(122, 12).Value = 3042
$ws.Cells.Item(122, 13).Value = -1746.6001
$ws.Cells.Item(126, 8).Value = 3713.8572
$ws.Cells.Item(126, 9).Value = 1999
$ws.Cells.Item(126, 10).Value = 3999.6667
$ws.Cells.Item(126, 11).Value = 5997
$ws.Cells.Item(126, 12).Value = 11999.0001
$ws.Cells.Item(126, 13).Value = -3527
$ws.Cells.Item(126, 14).Value = -16939.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 138239
$ws.Cells.Item(2, 9).Value = 100041.18
$ws.Cells.Item(2, 10).Value = 222274.2
$ws.Cells.Item(2, 11).Value = 600247.08
$ws.Cells.Item(2, 12).Value = 1333645.2
$ws.Cells.Item(2, 13).Value = -600134.08
$ws.Cells.Item(2, 14).Value = -1333871.2
$ws.Cells.Item(11, 8).Value = 25017724
$ws.Cells.Item(11, 9).Value = 31272092
$ws.Cells.Item(11, 10).Value = 250
$ws.Cells.Item(11, 11).Value = 93816276
$ws.Cells.Item(11, 12).Value = 750
$ws.Cells.Item(11, 13).Value = -93816136
$ws.Cells.Item(11, 14).Value = -1030
$ws.Cells.Item(68, 8).Value = 7996.0835
$ws.Cells.Item(68, 9).Value = 1950
$ws.Cells.Item(68, 10).Value = 8545.727999999999
$ws.Cells.Item(68, 11).Value = 5850
$ws.Cells.Item(68, 12).Value = 25637.184
$ws.Cells.Item(68, 13).Value = -5039
$ws.Cells.Item(68, 14).Value = -27259.184
$ws.Cells.Item(71, 8).Value = 7996.0835
$ws.Cells.Item(71, 9).Value = 1950
$ws.Cells.Item(71, 10).Value = 8545.727999999999
$ws.Cells.Item(71, 11).Value = 17550
$ws.Cells.Item(71, 12).Value = 76911.552
$ws.Cells.Item(71, 13).Value = -13494
$ws.Cells.Item(71, 14).Value = -85023.552
$ws.Cells.Item(132, 8).Value = 3683.4
$ws.Cells.Item(132, 9).Value = 2758.5
$ws.Cells.Item(132, 10).Value = 4300
$ws.Cells.Item(132, 11).Value = 24826.5
$ws.Cells.Item(132, 12).Value = 38700
$ws.Cells.Item(132, 13).Value = -22296.5
$ws.Cells.Item(132, 14).Value = -43760

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1470.8334
$ws.Cells.Item(102, 9).Value = 1322.6
$ws.Cells.Item(102, 10).Value = 2212
$ws.Cells.Item(102, 11).Value = 1322.6
$ws.Cells.Item(102, 12).Value = 2212
$ws.Cells.Item(102, 13).Value = 299.4000000000001
$ws.Cells.Item(122, 8).Value = 1712.7142
$ws.Cells.Item(122, 9).Value = 1398.8
$ws.Cells.Item(122, 10).Value = 2497.5
$ws.Cells.Item(122, 11).Value = 4196.4
$ws.Cells.Item(122, 12).Value = 7492.5
$ws.Cells.Item(122, 13).Value = -1746.4
$ws.Cells.Item(122, 14).Value = -12392.5
$ws.Cells.Item(126, 8).Value = 2058.8
$ws.Cells.Item(126, 9).Value = 1531.5
$ws.Cells.Item(126, 10).Value = 2849.75
$ws.Cells.Item(126, 11).Value = 4594.5
$ws.Cells.Item(126, 12).Value = 8549.25
$ws.Cells.Item(126, 13).Value = -2124.5
$ws.Cells.Item(126, 14).Value = -13489.25
$ws.Cells.Item(132, 8).Value = 2816.3333
$ws.Cells.Item(132, 9).Value = 2816.3333
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 8448.999899999999
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -5918.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6564.5
$ws.Cells.Item(7, 9).Value = 3216.125
$ws.Cells.Item(7, 10).Value = 8477.857
$ws.Cells.Item(7, 11).Value = 3216.125
$ws.Cells.Item(7, 12).Value = 8477.857
$ws.Cells.Item(7, 13).Value = -3104.125
$ws.Cells.Item(40, 8).Value = 2222.1667
$ws.Cells.Item(40, 9).Value = 1866.7
$ws.Cells.Item(40, 10).Value = 3999.5
$ws.Cells.Item(40, 11).Value = 1866.7
$ws.Cells.Item(40, 12).Value = 3999.5
$ws.Cells.Item(40, 13).Value = -1730.7
$ws.Cells.Item(40, 14).Value = -4271.5
$ws.Cells.Item(93, 8).Value = 1961.1666
$ws.Cells.Item(93, 9).Value = 1787.75
$ws.Cells.Item(93, 10).Value = 2099.9
$ws.Cells.Item(93, 11).Value = 1787.75
$ws.Cells.Item(93, 12).Value = 2099.9
$ws.Cells.Item(93, 13).Value = -539.75
$ws.Cells.Item(93, 14).Value = -4595.9
$ws.Cells.Item(122, 8).Value = 5961.269
$ws.Cells.Item(122, 9).Value = 9217.875
$ws.Cells.Item(122, 10).Value = 4513.8887
$ws.Cells.Item(122, 11).Value = 27653.625
$ws.Cells.Item(122, 12).Value = 13541.6661
$ws.Cells.Item(122, 13).Value = -25203.625
$ws.Cells.Item(122, 14).Value = -18441.6661
$ws.Cells.Item(126, 8).Value = 6564.5
$ws.Cells.Item(126, 9).Value = 3216.125
$ws.Cells.Item(126, 10).Value = 8477.857
$ws.Cells.Item(126, 11).Value = 9648.375
$ws.Cells.Item(126, 12).Value = 25433.571
$ws.Cells.Item(126, 13).Value = -7178.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 825.2
$ws.Cells.Item(122, 9).Value = 881.625
$ws.Cells.Item(122, 10).Value = 599.5
$ws.Cells.Item(122, 11).Value = 2644.875
$ws.Cells.Item(122, 12).Value = 1798.5
$ws.Cells.Item(122, 13).Value = -194.875
$ws.Cells.Item(122, 14).Value = -6698.5
$ws.Cells.Item(126, 8).Value = 3993.4546
$ws.Cells.Item(126, 9).Value = 4004.5
$ws.Cells.Item(126, 10).Value = 3883
$ws.Cells.Item(126, 11).Value = 12013.5
$ws.Cells.Item(126, 12).Value = 11649
$ws.Cells.Item(126, 13).Value = -9543.5
$ws.Cells.Item(126, 14).Value = -16589
$ws.Cells.Item(132, 8).Value = 3095.2273
$ws.Cells.Item(132, 9).Value = 3277.8333
$ws.Cells.Item(132, 10).Value = 2273.5
$ws.Cells.Item(132, 11).Value = 9833.499899999999
$ws.Cells.Item(132, 12).Value = 6820.5
$ws.Cells.Item(132, 13).Value = -7303.499899999999
$ws.Cells.Item(132, 14).Value = -11880.5
